$d = $word.ActiveDocument

# 1. Header date on the first page: "2023" -> "2024" (year of the contract)
$null = $d.Content.Find.Execute("23 року", $true, $false, $false, $false, $false, $true, 1, $false, "24 року", 2)

# 2. Power-of-attorney clause: old order number/date -> new order number/date,
#    and the signer's authority is now acting rector ("в.о. ректора")
$null = $d.Content.Find.Execute("Доручення ректора № 65 від 01.05.2023р., ", $true, $false, $false, $false, $false, $true, 1, $false, "Доручення в.о. ректора № 70 від 09.07.2024р., ", 2)

# 3. Contract date reference further down: "2023р" -> "2024р"
$null = $d.Content.Find.Execute("____________ 2023р", $true, $false, $false, $false, $false, $true, 1, $false, "____________ 2024р", 2)

# 4. Total price figure: 179 600,00 -> 199 600,00 (keep the non-breaking space)
$null = $d.Content.Find.Execute("179 600,00", $true, $false, $false, $false, $false, $true, 1, $false, "199 600,00", 2)

# 5. Total price spelled out in words: "сто сімдесят дев'ять" -> "сто дев'яносто дев'ять"
$null = $d.Content.Find.Execute("сто сімдесят дев’ять тисяч шістсот", $true, $false, $false, $false, $false, $true, 1, $false, "сто дев’яносто дев’ять тисяч шістсот", 2)
